$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bulk update: column C ("Förändrad") for rows 2-498 changes from 45181 to 45182.
$ws.Range("C2:C498").Value = 45182

# Rows 495, 496, 498 had their "Beteckning" (A) and "Area (ha)" (G) values
# cyclically rotated: 495 -> gets old 496 values, 496 -> gets old 498 values,
# 498 -> gets old 495 values. Row 497 is unchanged apart from column C above.
$a495 = $ws.Cells.Item(495, 1).Value()
$g495 = $ws.Cells.Item(495, 7).Value()
$a496 = $ws.Cells.Item(496, 1).Value()
$g496 = $ws.Cells.Item(496, 7).Value()
$a498 = $ws.Cells.Item(498, 1).Value()
$g498 = $ws.Cells.Item(498, 7).Value()

$ws.Cells.Item(495, 1).Value = $a496
$ws.Cells.Item(495, 7).Value = $g496

$ws.Cells.Item(496, 1).Value = $a498
$ws.Cells.Item(496, 7).Value = $g498

$ws.Cells.Item(498, 1).Value = $a495
$ws.Cells.Item(498, 7).Value = $g495
